$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 49, shifting existing rows 49..121 down to 50..122.
$ws.Rows("49").Insert()

# Populate the newly inserted row 49 with the new weekly record (same
# Mercado / Region / Categoria / etc. as the surrounding rows, new date +
# volume + price figures).
$ws.Range("A49").Value = 3
$ws.Range("B49").Value = "Femacal de La Calera"
$ws.Range("C49").Value = "Coquimbo"
$ws.Range("D49").Value = 44579
$ws.Range("E49").Value = 5
$ws.Range("F49").Value = 100112052
$ws.Range("G49").Value = "Albahaca"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 85
$ws.Range("K49").Value = 4500
$ws.Range("L49").Value = 5000
$ws.Range("M49").Value = 4735
$ws.Range("N49").Value = "$/docena de matas"
$ws.Range("O49").Value = "Provincia de Quillota"
$ws.Range("P49").Value = 789
$ws.Range("Q49").Value = 6
$ws.Range("R49").Value = "Hortaliza"
